$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.842.51"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.38%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.993.00"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.72%  "
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.05"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.81%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.27"
$ws.Range("D6").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.988.99"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.77%  "
$ws.Range("E8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("E9").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.33%  "
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.18"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("E11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.67%  "
$ws.Range("E12").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.63%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.79"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.09%  "
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.118"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.490.36"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.59%  "
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.925.30"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.987.98"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.84%  "
$ws.Range("E18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.15"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.92%  "
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "433.29"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.95%  "
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.09"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.91%  "
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.660"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -7.29%  "
$ws.Range("E22").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.07%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.85"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.19%  "
$ws.Range("E24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.69"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.59%  "
$ws.Range("E25").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E27").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.46%  "
$ws.Range("E28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.17"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -8.16%  "
$ws.Range("E29").ClearFormats()

$ws.Range("B30").Value = "ImmutableX"

$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.87"
$ws.Range("D30").ClearFormats()

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.67%  "
$ws.Range("E30").ClearFormats()

$ws.Range("B31").Value = "EthereumClassic"

$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.39"
$ws.Range("D31").ClearFormats()

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.01%  "
$ws.Range("E31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.04"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -11.27%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0932"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.86%  "
$ws.Range("E33").ClearFormats()

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.70%  "
$ws.Range("E34").ClearFormats()

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.953"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.44%  "
$ws.Range("E35").ClearFormats()

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.57"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("E36").ClearFormats()

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.01"
$ws.Range("D37").ClearFormats()

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0664"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.63%  "
$ws.Range("E38").ClearFormats()

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.05%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.78"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("E40").ClearFormats()

$ws.Range("B41").Value = "Bittensor"

$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "374.78"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.46%  "
$ws.Range("E41").ClearFormats()

$ws.Range("B42").Value = "Kaspa"

$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.107"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.14%  "
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.681.91"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.78%  "
$ws.Range("E43").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.86%  "
$ws.Range("E44").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E45").ClearFormats()

$ws.Range("B46").Value = "Monero"

$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "120.81"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("E46").ClearFormats()

$ws.Range("B47").Value = "Arweave"

$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.95"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("E47").ClearFormats()

$ws.Range("B48").Value = "TheGraph"

$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.234"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.47%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.96"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.17%  "
$ws.Range("E49").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("E50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.26"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.30%  "
$ws.Range("E51").ClearFormats()
